$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels in row 6
$ws.Range("B6").Value = "qbar"
$ws.Range("C6").Value = "pk"
$ws.Range("D6").Value = "Delta p k"

# Update formulas to reflect Laplace/Taubin smoothing sign convention
$ws.Range("D7").Formula = "=B7-C7"
$ws.Range("D8:D17").Formula = "=B8-C8"
$ws.Range("C8").Formula = "=C7+E7"
$ws.Range("C9:C17").Formula = "=C8+E8"

# Update the active selection to match the new view state
$ws.Range("F9").Select()
